$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value (serial 45171 = 2023-09-02).
# Bump it by one day (to 45172 = 2023-09-03) for every data row (2..153).
$lastRow = 153
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
